$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.20%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11.55%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.35%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05732"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.86%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8599"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.30%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8812"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.97%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1367"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.93%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07082"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.28%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02867"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.66%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09396"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.18%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001518"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.00%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04141"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.00%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005957"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.53%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.500"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.62%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.070"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.63%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.180"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-5.51%"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3184"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.29%"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03271"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.44%"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1308"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.27%"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.535"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.55%"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1379"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.41%"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "One"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0005997"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.15%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001214"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.83%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004504"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.44%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "23.48%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001382"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-0.05%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03791"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.85%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005598"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-7.12%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1073"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.93%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002589"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "12.62%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009996"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "22.41%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005098"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.55%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.01%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08895"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-18.37%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "4.90%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
